# Append the latest batch of status-report rows to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data rows to add (Date, Hours, Comment), continuing after the last
# existing row (row 37).
$rows = @(
    @("2/13/2010", 1.5, "Researched sterilizable motor"),
    @("2/14/2010", 0.5, "Encoder speed requirement calculation"),
    @("2/14/2010", 2,   "Researched QNX installation methods and issues"),
    @("2/15/2010", 3,   "Ordered control system components"),
    @("2/18/2010", 1,   "Group Meeting on Skype"),
    @("2/18/2010", 0.5, "Gravity Compensation motor experiment")
)

$startRow = 38
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $dateCell = $ws.Cells.Item($r, 1)
    # Force the date-shaped text to be stored as a literal string instead of
    # being auto-recognised as a date serial: format as Text, type it in,
    # then drop the Text number-format again so the cell keeps using the
    # sheet's normal (default) style.
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $rows[$i][0]
    $dateCell.ClearFormats()
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$lastRow = $startRow + $rows.Count - 1
[void]$ws.Range("A$($lastRow + 1)").Select()
